# Update "想去人数" (want-to-go count) values in column F across all
# four worksheets, matching the regenerated gh-pages data snapshot.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 90
$ws.Range("F5").Value = 1674
$ws.Range("F7").Value = 3591
$ws.Range("F13").Value = 2074
$ws.Range("F14").Value = 628140
$ws.Range("F17").Value = 1285
$ws.Range("F20").Value = 1207
$ws.Range("F21").Value = 2009
$ws.Range("F24").Value = 1481
$ws.Range("F25").Value = 674
$ws.Range("F26").Value = 1444
$ws.Range("F29").Value = 1044
$ws.Range("F31").Value = 1047
$ws.Range("F34").Value = 1948
$ws.Range("F39").Value = 184
$ws.Range("F41").Value = 2451
$ws.Range("F44").Value = 2989

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F9").Value = 446
$ws.Range("F11").Value = 143611
$ws.Range("F12").Value = 143611
$ws.Range("F19").Value = 321
$ws.Range("F23").Value = 69
$ws.Range("F24").Value = 76
$ws.Range("F27").Value = 418
$ws.Range("F31").Value = 54
$ws.Range("F33").Value = 249

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 3085
$ws.Range("F8").Value = 781
$ws.Range("F9").Value = 1087
$ws.Range("F10").Value = 608
$ws.Range("F11").Value = 1525
$ws.Range("F12").Value = 454
$ws.Range("F13").Value = 1651

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 781
$ws.Range("F3").Value = 1087
$ws.Range("F4").Value = 608
$ws.Range("F6").Value = 1525
$ws.Range("F7").Value = 454
$ws.Range("F9").Value = 90
$ws.Range("F10").Value = 1651
$ws.Range("F11").Value = 3591
$ws.Range("F16").Value = 2074
$ws.Range("F17").Value = 628140
$ws.Range("F18").Value = 446
$ws.Range("F21").Value = 143611
$ws.Range("F22").Value = 1285
$ws.Range("F25").Value = 1207
$ws.Range("F26").Value = 2009
$ws.Range("F29").Value = 1481
$ws.Range("F30").Value = 674
$ws.Range("F32").Value = 1444
$ws.Range("F35").Value = 69
$ws.Range("F36").Value = 1044
$ws.Range("F37").Value = 1047
$ws.Range("F38").Value = 76
$ws.Range("F41").Value = 1948
$ws.Range("F47").Value = 2451
$ws.Range("F50").Value = 2989
